$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G7"   = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
    "G17"  = "ola.m.abdelfattah@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
    "G24"  = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
    "G34"  = "ola.m.abdelfattah@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
    "G35"  = "abdullah.elagrody@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg"
    "G45"  = "mohamed.saleem@med.asu.edu.eg, backup@backdoor.com, System, Rania.a.youssef@med.asu.edu.eg"
    "G48"  = "Aya.alaa-eldein@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
    "G51"  = "Monica.Eshak@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg"
    "G62"  = "mohamed.saleem@med.asu.edu.eg, backup@backdoor.com, System, Rania.a.youssef@med.asu.edu.eg"
    "G65"  = "Aya.alaa-eldein@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
    "G68"  = "Monica.Eshak@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg"
    "G72"  = "Omnia.Mohammed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
    "G76"  = "mohamed.saleem@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
    "G83"  = "afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg, Youstina.ibrahim@med.asu.edu.eg"
    "G85"  = "yasmin.m.senosy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
    "G98"  = "user@user.com, Walaa.h.ghanima@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg"
    "G99"  = "Walaa.h.ghanima@med.asu.edu.eg, user@user.com"
    "G100" = "Aya.alaa-eldein@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg"
    "G102" = "yasmin.m.senosy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
    "G105" = "System, backup@backdoor.com"
    "G109" = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
    "G116" = "enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
    "G119" = "yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
    "G122" = "System, backup@backdoor.com"
    "G126" = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
    "G133" = "enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
    "G136" = "yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg"
    "G143" = "afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
    "G149" = "Walaa.h.ghanima@med.asu.edu.eg, user@user.com"
    "G150" = "afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg, Youstina.ibrahim@med.asu.edu.eg"
    "G151" = "randa.rabea@med.asu.edu.eg, System, marian.samir@med.asu.edu.eg"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
